# Insert a new weekly price record at row 104 of the "Arándano (blue)"
# sheet (Macroferia Regional de Talca). All the rows that previously sat
# at 104-133 shift down to 105-134, and the new row 104 is a duplicate of
# the (old) row 104 record except its date (column D) is bumped forward
# from 2023-03-20 (45005) to 2023-03-30 (45015).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 104:133 down to 105:134, leaving a blank row 104 behind.
$ws.Rows.Item(104).Insert()

# Seed the new blank row with a copy of what is now row 105 (the old
# row 104 contents), then correct the date for the new entry.
$ws.Range("A105:T105").Copy()
$ws.Range("A104").PasteSpecial()
$ws.Range("D104").Value = 45015
